$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 63.04999923706055, 56.22000122070312, 63.43999862670898, 55.11000061035156, 134262856),
    @(3, 59.34000015258789, 60.11999893188477, 63.47000122070312, 58.9900016784668, 134262856),
    @(4, 60.45000076293945, 62, 62.16999816894531, 56.54999923706055, 134262856),
    @(5, 64, 65.95999908447266, 71.41000366210938, 61.2400016784668, 134262856),
    @(6, 59.93999862670898, 62.5, 63.36000061035156, 56.34000015258789, 134262856),
    @(7, 72.76000213623047, 76.69999694824219, 79.15000152587891, 71.29000091552734, 134262856),
    @(8, 71.34999847412109, 73.83000183105469, 75.23000335693359, 68.83999633789062, 134262856),
    @(9, 74.40000152587891, 70.66000366210938, 74.91000366210938, 70.25, 134262856),
    @(10, 66.33999633789062, 64.91000366210938, 67.33000183105469, 63, 134262856),
    @(11, 64.04000091552734, 64.02999877929688, 66.06999969482422, 61.84999847412109, 134262856),
    @(12, 66.25, 66.08999633789062, 67.94000244140625, 62.91999816894531, 134262856),
    @(13, 64.41999816894531, 65.09999847412109, 67.51999664306641, 64.05000305175781, 134262856),
    @(14, 74.86000061035156, 75.62000274658203, 80.26000213623047, 73.55000305175781, 134262856),
    @(15, 65.5, 63.97999954223633, 68.12000274658203, 62.02999877929688, 134262856),
    @(16, 65.20999908447266, 70.33999633789062, 70.76000213623047, 63.5099983215332, 134262856),
    @(17, 64.25, 63.7599983215332, 64.48000335693359, 59.2599983215332, 134262856),
    @(18, 58.11999893188477, 62.68999862670898, 62.79999923706055, 56.63999938964844, 134262856),
    @(19, 78.04000091552734, 90.81999969482422, 90.93000030517578, 77.54000091552734, 134262856),
    @(20, 81.26999664306641, 81.58999633789062, 84.62000274658203, 79.79000091552734, 134262856),
    @(21, 75.55999755859375, 68.83999633789062, 76.31999969482422, 67.76000213623047, 134262856),
    @(22, 73.48999786376953, 63.65000152587891, 74.05000305175781, 63.5, 134262856),
    @(23, 34.22999954223633, 42.40000152587891, 46.18000030517578, 31.61000061035156, 134262856),
    @(24, 38.18000030517578, 37.95999908447266, 39.43000030517578, 35.72999954223633, 134262856),
    @(25, 34.40000152587891, 33.90000152587891, 36.54999923706055, 32.88000106811523, 134262856),
    @(26, 32.70000076293945, 43.54000091552734, 46.79999923706055, 31.15999984741211, 134262856),
    @(27, 50.47000122070312, 51.79000091552734, 53.34999847412109, 46.29999923706055, 134262856),
    @(28, 49.9900016784668, 49.63999938964844, 51.86999893188477, 46.61000061035156, 134262856),
    @(29, 55.56999969482422, 59.68999862670898, 59.9900016784668, 53.77000045776367, 134262856),
    @(30, 44.40000152587891, 44.02000045776367, 47.90000152587891, 43.47000122070312, 134262856),
    @(31, 48.79000091552734, 36.81000137329102, 50, 36.20000076293945, 134262856),
    @(32, 30.05999946594238, 32.93000030517578, 33.18000030517578, 25.3799991607666, 134262856),
    @(33, 35.56000137329102, 40.95999908447266, 44.77000045776367, 35.0099983215332, 134262856),
    @(34, 31.89999961853028, 34.45000076293945, 38.2599983215332, 31.67000007629395, 134262856),
    @(35, 33.79999923706055, 35.02999877929688, 36.25, 32.02000045776367, 134262856),
    @(36, 41.11000061035156, 30.94000053405762, 44.31999969482422, 27.34000015258789, 134262856),
    @(37, 18.3700008392334, 18.44000053405762, 22.21999931335449, 15.02000045776367, 134262856),
    @(38, 27.64999961853028, 22.22999954223633, 29.11000061035156, 22.17000007629395, 134262856),
    @(39, 18.09000015258789, 15.90999984741211, 18.28000068664551, 15.02000045776367, 134262856),
    @(40, 12.64000034332275, 20.21999931335449, 21.31999969482422, 12.30000019073486, 134262856),
    @(41, 11.90999984741211, 9.600000381469728, 12.01000022888184, 9.529999732971191, 134262856),
    @(42, 8.800000190734863, 9.619999885559082, 12.55000019073486, 7.449999809265137, 134262856),
    @(43, 10.47999954223633, 9.270000457763672, 10.48999977111816, 7.360000133514404, 134262856),
    @(44, 14.65999984741211, 16.43000030517578, 17.26000022888184, 14, 134262856)
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 4).Value = $item[1]
    $ws.Cells.Item($r, 5).Value = $item[2]
    $ws.Cells.Item($r, 6).Value = $item[3]
    $ws.Cells.Item($r, 7).Value = $item[4]
    $ws.Cells.Item($r, 8).Value = $item[5]
    $ws.Cells.Item($r, 9).Value = "VSAT"
}